$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = "L"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = "m2"
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = "W/m2K"
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = "°C"
